{"js": "// Each entry maps the original 3-digit-by-1-digit multiplication problem\n// (as it appears in the table cell) to its replacement.\nconst replacements = [\n  [\"173\u00d79=1557\", \"430\u00d72=860\"],\n  [\"718\u00d75=3590\", \"943\u00d72=1886\"],\n  [\"271\u00d73=813\", \"647\u00d78=5176\"],\n  [\"497\u00d77=3479\", \"741\u00d75=3705\"],\n  [\"557\u00d79=5013\", \"352\u00d75=1760\"],\n  [\"788\u00d77=5516\", \"617\u00d73=1851\"],\n  [\"826\u00d78=6608\", \"598\u00d73=1794\"],\n  [\"338\u00d79=3042\", \"314\u00d72=628\"],\n  [\"397\u00d76=2382\", \"504\u00d73=1512\"],\n  [\"776\u00d74=3104\", \"144\u00d72=288\"],\n  [\"637\u00d76=3822\", \"861\u00d73=2583\"],\n  [\"521\u00d76=3126\", \"551\u00d73=1653\"],\n  [\"875\u00d76=5250\", \"696\u00d75=3480\"],\n  [\"952\u00d74=3808\", \"942\u00d75=4710\"],\n  [\"849\u00d74=3396\", \"388\u00d78=3104\"],\n  [\"938\u00d76=5628\", \"945\u00d76=5670\"],\n  [\"305\u00d78=2440\", \"159\u00d73=477\"],\n  [\"500\u00d78=4000\", \"760\u00d76=4560\"],\n  [\"338\u00d73=1014\", \"242\u00d74=968\"],\n  [\"997\u00d75=4985\", \"769\u00d73=2307\"],\n  [\"449\u00d79=4041\", \"191\u00d78=1528\"],\n  [\"412\u00d78=3296\", \"811\u00d78=6488\"],\n  [\"165\u00d73=495\", \"777\u00d73=2331\"],\n  [\"636\u00d76=3816\", \"705\u00d72=1410\"],\n  [\"400\u00d78=3200\", \"534\u00d77=3738\"],\n];\n\nfor (const [before, after] of replacements) {\n  // Find the table cell containing the old equation text.\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${before}\"`);\n  }\n\n  // Replace the full run's text (there is exactly one match per value).\n  for (const range of results.items) {\n    range.insertText(after, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each old 3-digit-by-1-digit multiplication equation in the table\n# with its new equation. Every value in the answer table is unique, so a\n# plain Find/Replace (wdReplaceAll) on each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('173\u00d79=1557', '430\u00d72=860'),\n    @('718\u00d75=3590', '943\u00d72=1886'),\n    @('271\u00d73=813', '647\u00d78=5176'),\n    @('497\u00d77=3479', '741\u00d75=3705'),\n    @('557\u00d79=5013', '352\u00d75=1760'),\n    @('788\u00d77=5516', '617\u00d73=1851'),\n    @('826\u00d78=6608', '598\u00d73=1794'),\n    @('338\u00d79=3042', '314\u00d72=628'),\n    @('397\u00d76=2382', '504\u00d73=1512'),\n    @('776\u00d74=3104', '144\u00d72=288'),\n    @('637\u00d76=3822', '861\u00d73=2583'),\n    @('521\u00d76=3126', '551\u00d73=1653'),\n    @('875\u00d76=5250', '696\u00d75=3480'),\n    @('952\u00d74=3808', '942\u00d75=4710'),\n    @('849\u00d74=3396', '388\u00d78=3104'),\n    @('938\u00d76=5628', '945\u00d76=5670'),\n    @('305\u00d78=2440', '159\u00d73=477'),\n    @('500\u00d78=4000', '760\u00d76=4560'),\n    @('338\u00d73=1014', '242\u00d74=968'),\n    @('997\u00d75=4985', '769\u00d73=2307'),\n    @('449\u00d79=4041', '191\u00d78=1528'),\n    @('412\u00d78=3296', '811\u00d78=6488'),\n    @('165\u00d73=495', '777\u00d73=2331'),\n    @('636\u00d76=3816', '705\u00d72=1410'),\n    @('400\u00d78=3200', '534\u00d77=3738'),\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute(\n        $findText, $true, $false, $false, $false, $false, $true, 1, $false,\n        $replaceText, 2     # 2 = wdReplaceAll\n    )\n\n    if (-not $found) {\n        throw \"No match found for '$findText'\"\n    }\n}\n"}
